$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4712233308538032
$ws.Range("C2").Value = 0.1274291188181138
$ws.Range("D2").Value = 0.04711151860685447
$ws.Range("E2").Value = 0.1138147368154083
$ws.Range("F2").Value = 0.9367715181862977
$ws.Range("I2").Value = 0.8793434615819109
$ws.Range("K2").Value = 0.2882797096711442
$ws.Range("L2").Value = 0.2091133763820139
$ws.Range("M2").Value = 0.1483359338671875
$ws.Range("N2").Value = 1.8697055941516
$ws.Range("O2").Value = 3.342655214366857
$ws.Range("B3").Value = 0.4369290482055703
$ws.Range("C3").Value = 0.1253961860420318
$ws.Range("D3").Value = 0.04480277232477192
$ws.Range("E3").Value = 0.1142244436021809
$ws.Range("F3").Value = 0.937333887060035
$ws.Range("I3").Value = 0.8851256903599563
$ws.Range("K3").Value = 0.2553132173878652
$ws.Range("L3").Value = 0.2064207288184434
$ws.Range("M3").Value = 0.1417215665227225
$ws.Range("N3").Value = 1.887282939824448
$ws.Range("O3").Value = 3.357312927166845
$ws.Range("B4").Value = 0.4160029602486759
$ws.Range("C4").Value = 0.1241315927486966
$ws.Range("D4").Value = 0.04336901101641644
$ws.Range("E4").Value = 0.1145169273562843
$ws.Range("F4").Value = 0.9381647070590091
$ws.Range("I4").Value = 0.8890516628248832
$ws.Range("K4").Value = 0.2350736876104662
$ws.Range("L4").Value = 0.2048670834967581
$ws.Range("M4").Value = 0.1377204289288336
$ws.Range("N4").Value = 1.898628565866457
$ws.Range("O4").Value = 3.367895021276709
$ws.Range("B5").Value = 0.4075088987898994
$ws.Range("C5").Value = 0.123612163087067
$ws.Range("D5").Value = 0.04278069542084495
$ws.Range("E5").Value = 0.1146464275102659
$ws.Range("F5").Value = 0.9386254845195126
$ws.Range("I5").Value = 0.8907460646974954
$ws.Range("K5").Value = 0.2268269068872542
$ws.Range("L5").Value = 0.2042590956540238
$ws.Range("M5").Value = 0.1361051756573453
$ws.Range("N5").Value = 1.903391114857171
$ws.Range("O5").Value = 3.372605473493451
$ws.Range("B6").Value = 0.4061005070761041
$ws.Range("C6").Value = 0.1235256652118863
$ws.Range("D6").Value = 0.04268276235237778
$ws.Range("E6").Value = 0.1146685542510895
$ws.Range("F6").Value = 0.9387093811539202
$ws.Range("I6").Value = 0.8910331307512145
$ws.Range("K6").Value = 0.2254576093132101
$ws.Range("L6").Value = 0.204159660109859
$ws.Range("M6").Value = 0.1358378884570826
$ws.Range("N6").Value = 1.904190335905513
$ws.Range("O6").Value = 3.37341169757353
$ws.Range("B7").Value = 0.4158882699551043
$ws.Range("C7").Value = 0.1241246041058233
$ws.Range("D7").Value = 0.04336109313921099
$ws.Range("E7").Value = 0.114518632065872
$ws.Range("F7").Value = 0.938170426272201
$ws.Range("I7").Value = 0.8890741312502719
$ws.Range("K7").Value = 0.2349624639871735
$ws.Range("L7").Value = 0.2048587820862906
$ws.Range("M7").Value = 0.1376985831718436
$ws.Range("N7").Value = 1.898692231943724
$ws.Range("O7").Value = 3.367956935530401
$ws.Range("B8").Value = 0.4593718841306327
$ws.Range("C8").Value = 0.126731569881926
$ws.Range("D8").Value = 0.04631883451543217
$ws.Range("E8").Value = 0.1139475224948843
$ws.Range("F8").Value = 0.9368647240680659
$ws.Range("I8").Value = 0.8812592493873446
$ws.Range("K8").Value = 0.2769127208750319
$ws.Range("L8").Value = 0.2081643138938603
$ws.Range("M8").Value = 0.1460428982162014
$ws.Range("N8").Value = 1.875651501021022
$ws.Range("O8").Value = 3.347380988038239
$ws.Range("B9").Value = 0.545657917297774
$ws.Range("C9").Value = 0.1317134242825162
$ws.Range("D9").Value = 0.05198979441599505
$ws.Range("E9").Value = 0.1131514229903594
$ws.Range("F9").Value = 0.9381523226505735
$ws.Range("I9").Value = 0.8689119607085303
$ws.Range("K9").Value = 0.3591759359457853
$ws.Range("L9").Value = 0.2154342973658885
$ws.Range("M9").Value = 0.1628783510069312
$ws.Range("N9").Value = 1.834854470795763
$ws.Range("O9").Value = 3.319574191233016
$ws.Range("B10").Value = 0.6096470512793815
$ws.Range("C10").Value = 0.1352936603265462
$ws.Range("D10").Value = 0.05607682412168913
$ws.Range("E10").Value = 0.1127628699070566
$ws.Range("F10").Value = 0.9414399595125289
$ws.Range("I10").Value = 0.8616518914159741
$ws.Range("K10").Value = 0.4195967128333393
$ws.Range("L10").Value = 0.2212528341921853
$ws.Range("M10").Value = 0.1755303499955545
$ws.Range("N10").Value = 1.807550091563199
$ws.Range("O10").Value = 3.306778344179207
$ws.Range("B11").Value = 0.638881828771332
$ws.Range("C11").Value = 0.1369049701394545
$ws.Range("D11").Value = 0.05791874631211869
$ws.Range("E11").Value = 0.1126285159728067
$ws.Range("F11").Value = 0.9434432528822896
$ws.Range("I11").Value = 0.8587416733104192
$ws.Range("K11").Value = 0.4470764329413441
$ws.Range("L11").Value = 0.2240028821448306
$ws.Range("M11").Value = 0.18134655543804
$ws.Range("N11").Value = 1.795707255514241
$ws.Range("O11").Value = 3.30261236735609
$ws.Range("B12").Value = 0.6499698249714356
$ws.Range("C12").Value = 0.1375126213511066
$ws.Range("D12").Value = 0.05861372882610993
$ws.Range("E12").Value = 0.1125837173328179
$ws.Range("F12").Value = 0.9442747617211822
$ws.Range("I12").Value = 0.8576960157322446
$ws.Range("K12").Value = 0.4574810075855567
$ws.Range("L12").Value = 0.2250590215760297
$ws.Range("M12").Value = 0.1835576242304668
$ws.Range("N12").Value = 1.791305769015727
$ws.Range("O12").Value = 3.301272551874206
$ws.Range("B13").Value = 0.6475810624865801
$ws.Range("C13").Value = 0.1373818650598935
$ws.Range("D13").Value = 0.05846416413204736
$ws.Range("E13").Value = 0.1125930954892951
$ws.Range("F13").Value = 0.9440924406369717
$ws.Range("I13").Value = 0.8579187101507983
$ws.Range("K13").Value = 0.4552402670532558
$ws.Range("L13").Value = 0.2248309080426907
$ws.Range("M13").Value = 0.183081050796666
$ws.Range("N13").Value = 1.792250010035843
$ws.Range("O13").Value = 3.301550534738794
$ws.Range("B14").Value = 0.6397936995971065
$ws.Range("C14").Value = 0.1369550125303647
$ws.Range("D14").Value = 0.05797597347240924
$ws.Range("E14").Value = 0.1126247086522412
$ws.Range("F14").Value = 0.9435102010977872
$ws.Range("I14").Value = 0.8586545168088477
$ws.Range("K14").Value = 0.4479324546755663
$ws.Range("L14").Value = 0.2240894760951591
$ws.Range("M14").Value = 0.1815282898692701
$ws.Range("N14").Value = 1.795343476265066
$ws.Range("O14").Value = 3.302497376370809
$ws.Range("B15").Value = 0.6350259637482623
$ws.Range("C15").Value = 0.1366932246451142
$ws.Range("D15").Value = 0.05767661445645444
$ws.Range("E15").Value = 0.1126448636378754
$ws.Range("F15").Value = 0.9431630538071545
$ws.Range("I15").Value = 0.8591125603061229
$ws.Range("K15").Value = 0.4434560104456295
$ws.Range("L15").Value = 0.2236372469231753
$ws.Range("M15").Value = 0.1805782945992007
$ws.Range("N15").Value = 1.797249141185607
$ws.Range("O15").Value = 3.303108299505595
$ws.Range("B16").Value = 0.6077389640202853
$ws.Range("C16").Value = 0.1351880063811564
$ws.Range("D16").Value = 0.05595609995077666
$ws.Range("E16").Value = 0.1127725018489585
$ws.Range("F16").Value = 0.9413192473837739
$ws.Range("I16").Value = 0.8618499732044214
$ws.Range("K16").Value = 0.4178006819924747
$ws.Range("L16").Value = 0.2210751812945517
$ws.Range("M16").Value = 0.1751514579485445
$ws.Range("N16").Value = 1.808335681976452
$ws.Range("O16").Value = 3.307083882470664
$ws.Range("B17").Value = 0.5910310177411873
$ws.Range("C17").Value = 0.1342601432230381
$ws.Range("D17").Value = 0.05489617265241264
$ws.Range("E17").Value = 0.1128616501654864
$ws.Range("F17").Value = 0.9403180987911099
$ws.Range("I17").Value = 0.8636297600435796
$ws.Range("K17").Value = 0.4020600518599338
$ws.Range("L17").Value = 0.2195298041020237
$ws.Range("M17").Value = 0.1718377306756693
$ws.Range("N17").Value = 1.81528498048764
$ws.Range("O17").Value = 3.309946477742869
$ws.Range("B18").Value = 0.5814329298522409
$ws.Range("C18").Value = 0.1337248283536923
$ws.Range("D18").Value = 0.05428490491298987
$ws.Range("E18").Value = 0.1129169180141911
$ws.Range("F18").Value = 0.9397900610194299
$ws.Range("I18").Value = 0.864690386202458
$ws.Range("K18").Value = 0.3930059256684331
$ws.Range("L18").Value = 0.2186506592768609
$ws.Range("M18").Value = 0.169937488480322
$ws.Range("N18").Value = 1.819336442241072
$ws.Range("O18").Value = 3.311748749093084
$ws.Range("B19").Value = 0.5781852431519496
$ws.Range("C19").Value = 0.1335433001834403
$ws.Range("D19").Value = 0.05407766179578744
$ws.Range("E19").Value = 0.1129363169805604
$ws.Range("F19").Value = 0.9396194887997709
$ws.Range("I19").Value = 0.8650558420121293
$ws.Range("K19").Value = 0.3899402790774786
$ws.Range("L19").Value = 0.2183546669602663
$ws.Range("M19").Value = 0.1692950867675194
$ws.Range("N19").Value = 1.82071754200634
$ws.Range("O19").Value = 3.312385729827781
$ws.Range("B20").Value = 0.5928083823627617
$ws.Range("C20").Value = 0.1343590849504039
$ws.Range("D20").Value = 0.05500917216343737
$ws.Range("E20").Value = 0.1128517471453083
$ws.Range("F20").Value = 0.9404197267826717
$ws.Range("I20").Value = 0.8634364759320476
$ws.Range("K20").Value = 0.4037357277624665
$ws.Range("L20").Value = 0.2196933072229683
$ws.Range("M20").Value = 0.1721898909670045
$ws.Range("N20").Value = 1.814539584571816
$ws.Range("O20").Value = 3.30962562834074
$ws.Range("B21").Value = 0.6420805693989848
$ws.Range("C21").Value = 0.1370804579386373
$ws.Range("D21").Value = 0.05811943533664987
$ws.Range("E21").Value = 0.1126152582879083
$ws.Range("F21").Value = 0.9436792411889456
$ws.Range("I21").Value = 0.8584368628589303
$ws.Range("K21").Value = 0.4500789785664665
$ws.Range("L21").Value = 0.2243068527910168
$ws.Range("M21").Value = 0.1819841411418892
$ws.Range("N21").Value = 1.794432593771076
$ws.Range("O21").Value = 3.302212815380273
$ws.Range("B22").Value = 0.6743840204542266
$ws.Range("C22").Value = 0.138844352705874
$ws.Range("D22").Value = 0.06013750901306025
$ws.Range("E22").Value = 0.1124961192462806
$ws.Range("F22").Value = 0.9462344193155161
$ws.Range("I22").Value = 0.8554979221247265
$ws.Range("K22").Value = 0.4803585678714057
$ws.Range("L22").Value = 0.2274080435502697
$ws.Range("M22").Value = 0.1884353016944758
$ws.Range("N22").Value = 1.781776115710125
$ws.Range("O22").Value = 3.298753794513942
$ws.Range("B23").Value = 0.6571340196217932
$ws.Range("C23").Value = 0.1379042789728118
$ws.Range("D23").Value = 0.05906177588741457
$ws.Range("E23").Value = 0.112556471296573
$ws.Range("F23").Value = 0.9448318297856488
$ws.Range("I23").Value = 0.8570364405006501
$ws.Range("K23").Value = 0.4641987278300519
$ws.Range("L23").Value = 0.2257450394847069
$ws.Range("M23").Value = 0.1849876616198216
$ws.Range("N23").Value = 1.788486769917317
$ws.Range("O23").Value = 3.300473223918459
$ws.Range("B24").Value = 0.5920048123023491
$ws.Range("C24").Value = 0.1343143592325617
$ws.Range("D24").Value = 0.05495809100536064
$ws.Range("E24").Value = 0.112856211792252
$ws.Range("F24").Value = 0.9403736326895071
$ws.Range("I24").Value = 0.8635237431673488
$ws.Range("K24").Value = 0.4029781690630614
$ws.Range("L24").Value = 0.2196193584468915
$ws.Range("M24").Value = 0.1720306641391218
$ws.Range("N24").Value = 1.814876402951517
$ws.Range("O24").Value = 3.309770196766976
$ws.Range("B25").Value = 0.5222090826934789
$ws.Range("C25").Value = 0.1303797058565408
$ws.Range("D25").Value = 0.05046954738889298
$ws.Range("E25").Value = 0.1133322312601938
$ws.Range("F25").Value = 0.9373924672954885
$ws.Range("I25").Value = 0.8719338412972277
$ws.Range("K25").Value = 0.3369234643374455
$ws.Range("L25").Value = 0.2133835404093389
$ws.Range("M25").Value = 0.1582738623623783
$ws.Range("N25").Value = 1.845422094532086
$ws.Range("O25").Value = 3.325755187718499
